$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2-31) down to rows 13-42 to make room for 11 new rows
$ws.Range("A2:E31").Cut($ws.Range("A1000"))
$ws.Range("A1000:E1029").Cut($ws.Range("A13"))
$ws.Range("A1000:E1029").Delete(-4162) | Out-Null

# Populate new rows 2-12 with backward-extended data (1985-1995)
$ws.Cells.Item(2,1).Value = 31228
$ws.Cells.Item(2,2).Value = 1985
$ws.Cells.Item(2,3).Value = 1.274495535013775
$ws.Cells.Item(2,4).Value = 1986
$ws.Cells.Item(2,5).Value = 0.6162543489843353
$ws.Cells.Item(3,1).Value = 31593
$ws.Cells.Item(3,2).Value = 1986
$ws.Cells.Item(3,3).Value = 0.7987811539740708
$ws.Cells.Item(3,4).Value = 1987
$ws.Cells.Item(3,5).Value = 0.07643019850551713
$ws.Cells.Item(4,1).Value = 31958
$ws.Cells.Item(4,2).Value = 1987
$ws.Cells.Item(4,3).Value = -2.275857798620984
$ws.Cells.Item(4,4).Value = 1988
$ws.Cells.Item(4,5).Value = -3.016613995130724
$ws.Cells.Item(5,1).Value = 32324
$ws.Cells.Item(5,2).Value = 1988
$ws.Cells.Item(5,3).Value = 1.704500121876951
$ws.Cells.Item(5,4).Value = 1989
$ws.Cells.Item(5,5).Value = 1.379850795239923
$ws.Cells.Item(6,1).Value = 32689
$ws.Cells.Item(6,2).Value = 1989
$ws.Cells.Item(6,3).Value = 4.674707056384286
$ws.Cells.Item(6,4).Value = 1990
$ws.Cells.Item(6,5).Value = 4.526405235077502
$ws.Cells.Item(7,1).Value = 33054
$ws.Cells.Item(7,2).Value = 1990
$ws.Cells.Item(7,3).Value = 6.08517480942079
$ws.Cells.Item(7,4).Value = 1991
$ws.Cells.Item(7,5).Value = 6.655798235769961
$ws.Cells.Item(8,1).Value = 33419
$ws.Cells.Item(8,2).Value = 1991
$ws.Cells.Item(8,3).Value = 8.961416689426937
$ws.Cells.Item(8,4).Value = 1992
$ws.Cells.Item(8,5).Value = 9.375632354513574
$ws.Cells.Item(9,1).Value = 33785
$ws.Cells.Item(9,2).Value = 1992
$ws.Cells.Item(9,3).Value = 3.692571671353462
$ws.Cells.Item(9,4).Value = 1993
$ws.Cells.Item(9,5).Value = 4.510931229133397
$ws.Cells.Item(10,1).Value = 34150
$ws.Cells.Item(10,2).Value = 1993
$ws.Cells.Item(10,3).Value = -3.382453462135548
$ws.Cells.Item(10,4).Value = 1994
$ws.Cells.Item(10,5).Value = -3.798381663081463
$ws.Cells.Item(11,1).Value = 34515
$ws.Cells.Item(11,2).Value = 1994
$ws.Cells.Item(11,3).Value = 2.305922482069911
$ws.Cells.Item(11,4).Value = 1995
$ws.Cells.Item(11,5).Value = 1.982627263179704
$ws.Cells.Item(12,1).Value = 34880
$ws.Cells.Item(12,2).Value = 1995
$ws.Cells.Item(12,3).Value = 1.930360200805081
$ws.Cells.Item(12,4).Value = 1996
$ws.Cells.Item(12,5).Value = 1.674662359455792